# Course-view fix-up per 5/3 meeting: shrink the EE467/EE311 row's
# over-tall auto height, and add the new EE106 course (bulk populator).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42 (EE467/EE311): it had auto-grown to Excel's max row height
#     (409.5pt); set it to an explicit, human-sized custom height instead. ---
$ws.Rows(42).RowHeight = 140.25

# --- New course row: EE106 (Electronic Design For Software Development 1) ---

$ws.Range("A43").Value = 'EE106'
$ws.Range("B43").Value = "Electronic Design For Software`nDevelopment 1"
$ws.Range("C43").Value = '20'
$ws.Range("D43").Value = '1/2'
$ws.Range("E43").Value = 'NE'
$ws.Range("F43").Value = '1'
$ws.Range("G43").Value = 'This module will provide an introduction to the practice of computer programming. It is assumed troughout that'
$ws.Range("H43").Value = 'LO 1: To demonstrate a deep understanding of primitive data types, then later abstract data types.'
$ws.Range("I43").Value = 'This module will begin by teaching how to use the in-built Python GUI (IDLE) and interpretor.'
$ws.Range("J43").Value = 'LO1'
$ws.Range("K43").Value = 'Students are provided with a score at the end of a multiple choice exam at the end of semester 1. Feedback is'
$ws.Range("L43").Value = 'Overall average is calculated as follows: Demonstrated Exercises (30%), Semester 1 test (10%), Interim Report'
$ws.Range("M43").Value = 'Dive into Python 3. Other free material available online.'
$ws.Range("O43").Value = '13'
$ws.Range("P43").Value = '0'
$ws.Range("Q43").Value = '0'
$ws.Range("R43").Value = '66'
$ws.Range("S43").Value = '121'
$ws.Range("B43").WrapText = $true
$ws.Rows(43).RowHeight = 45

$ws.Range("G44").Value = 'the students have no prior exposure to programming.'
$ws.Range("H44").Value = 'LO 2: To demonstrate a deep understanding of basic programming constructions: if-else, for, while etc.'
$ws.Range("I44").Value = 'The analogy of a computer program to a cooking recipe is used teach the dsitinction between data types and'
$ws.Range("J44").Value = 'C 1: Progression through work sets'
$ws.Range("K44").Value = 'also provided after assessment of an interim report on week 5 of semester 2.'
$ws.Range("L44").Value = '(20%), Final Report (20%), Demonstraion (20%).'

$ws.Range("G45").Value = 'The course begins with an introduction to data types and their associated operators, before examinging basic'
$ws.Range("H45").Value = 'LO 3: To develop an understanding of systematic problem solving via the divide-and-conquer approach of'
$ws.Range("I45").Value = 'algorithms.'
$ws.Range("J45").Value = 'C 2: Semester 1 multiple choice exam'

$ws.Range("G46").Value = 'flow control in the form of decisions (if-else) and iteration (for/while loops). Problem decomposition through'
$ws.Range("H46").Value = 'functional decomposition.'
$ws.Range("I46").Value = 'The first semester is run as taught course in 3 main subsections: Data Types & Operators, Algorithmic'
$ws.Range("J46").Value = 'LO2'

$ws.Range("G47").Value = 'functions is discussed. Bespoke data types are addressed, and finally object oriented design principles are'
$ws.Range("H47").Value = 'LO 4: To develop self-reliance on programming by finding solutions to problems on own initiative.'
$ws.Range("I47").Value = 'Components (if-else, for, while, etc.), and program structure (functions, objects). In this semester, various'
$ws.Range("J47").Value = 'C 1: Progression through work sets'

$ws.Range("G48").Value = 'taught (inheritance and polymorphism).'
$ws.Range("H48").Value = 'LO 5: To learn to document code and produce software reports'
$ws.Range("I48").Value = 'problem sets on these concepts are provided.'
$ws.Range("J48").Value = 'C 2: Semester 1 exam'

$ws.Range("G49").Value = 'The first semester is based upon small self-contained problem sets to encourange problem solving skills while'
$ws.Range("I49").Value = 'The second semester is run as a programming project section. Groups of 2 are assigned a topic that is to be'
$ws.Range("J49").Value = 'LO3'

$ws.Range("G50").Value = 'cementing core concepts. The second semester involves small focussed projects undertaken in pairs.'
$ws.Range("I50").Value = 'completed over a 10 week period. The emphasis in this semester is on problem solving and self-reliance.'
$ws.Range("J50").Value = 'C 1: Second semester software design project: interim report, final report, demonstration'

$ws.Range("G51").Value = 'The programming language used in this module is Python; this permits the programmer to concentrate on the'
$ws.Range("J51").Value = 'LO4'

$ws.Range("G52").Value = 'problem solving aspects and core techniques of programming rather than being distracted by the syntax of the'
$ws.Range("J52").Value = 'C 1: Second semester software design project: interim report, final report, demonstration'

$ws.Range("G53").Value = 'language.'

# --- Leave the workbook scrolled/selected the way it was left in the
#     session (cursor on the row after the new EE106 entry). ---
$ws.Range("A44").Select()
